# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-job leve sheets. Some rows lose their
# LeveProfit (M/N) cells entirely when the underlying price feed no longer
# returns a value for that combination; two CRP rows gain a fresh
# LeveProfitHQ (N) cell where none existed before.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 61998.54
$ws.Range("I11").Value = 61998.54
$ws.Range("K11").Value = 61998.54
$ws.Range("M11").Value = -61858.54
$ws.Range("H19").Value = 1174.5
$ws.Range("J19").Value = 1174.5
$ws.Range("L19").Value = 1174.5
$ws.Range("N19").Value = -1524.5
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H100").Value = 34257.84
$ws.Range("I100").Value = 37651.57
$ws.Range("K100").Value = 37651.57
$ws.Range("M100").Value = -37110.57
$ws.Range("H113").Value = 3875
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3875
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3875
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10383
$ws.Range("H138").Value = 2506.919
$ws.Range("J138").Value = 3142.647
$ws.Range("L138").Value = 9427.940999999999
$ws.Range("N138").Value = -19707.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5750
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 5750
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 5750
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -6038
$ws.Range("H32").Value = 6878.125
$ws.Range("I32").Value = 4634.6
$ws.Range("K32").Value = 4634.6
$ws.Range("M32").Value = -4347.6
$ws.Range("H45").Value = 10411.462
$ws.Range("I45").Value = 11029.083
$ws.Range("K45").Value = 11029.083
$ws.Range("M45").Value = -10652.083
$ws.Range("H97").Value = 642.61365
$ws.Range("I97").Value = 634.175
$ws.Range("J97").Value = 727
$ws.Range("K97").Value = 634.175
$ws.Range("L97").Value = 727
$ws.Range("M97").Value = -138.175
$ws.Range("N97").Value = -1719
$ws.Range("H102").Value = 4635.684
$ws.Range("I102").Value = 4671
$ws.Range("K102").Value = 4671
$ws.Range("M102").Value = -3049

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 434.65
$ws.Range("I94").Value = 422.42105
$ws.Range("K94").Value = 422.42105
$ws.Range("M94").Value = 28.57895000000002
$ws.Range("H134").Value = 10523.6
$ws.Range("I134").Value = 8862.5
$ws.Range("J134").Value = 12184.7
$ws.Range("K134").Value = 26587.5
$ws.Range("L134").Value = 36554.10000000001
$ws.Range("M134").Value = -24052.5
$ws.Range("N134").Value = -41624.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 210.18182
$ws.Range("J7").Value = 521.5
$ws.Range("L7").Value = 521.5
$ws.Range("N7").Value = -747.5
$ws.Range("H13").Value = 15562.375
$ws.Range("J13").Value = 19499
$ws.Range("L13").Value = 19499
$ws.Range("N13").Value = -19777
$ws.Range("H16").Value = 2710.4666
$ws.Range("I16").Value = 2606.3076
$ws.Range("K16").Value = 2606.3076
$ws.Range("M16").Value = -2319.3076
$ws.Range("H22").Value = 323.5
$ws.Range("J22").Value = 337.6
$ws.Range("L22").Value = 337.6
$ws.Range("N22").Value = -1037.6
$ws.Range("H31").Value = 2654.6072
$ws.Range("J31").Value = 4457.9
$ws.Range("L31").Value = 4457.9
$ws.Range("N31").Value = -5047.9
$ws.Range("H34").Value = 2654.6072
$ws.Range("J34").Value = 4457.9
$ws.Range("L34").Value = 4457.9
$ws.Range("N34").Value = -4861.9
$ws.Range("H58").Value = 5176.5386
$ws.Range("I58").Value = 2548.75
$ws.Range("K58").Value = 2548.75
$ws.Range("M58").Value = -2345.75
$ws.Range("H62").Value = 501503
$ws.Range("J62").Value = 3006
$ws.Range("L62").Value = 3006
$ws.Range("N62").Value = -4254
$ws.Range("H65").Value = 501503
$ws.Range("J65").Value = 3006
$ws.Range("L65").Value = 15030
$ws.Range("N65").Value = -21270
$ws.Range("H74").Value = 42474.5
$ws.Range("J74").Value = 42474.5
$ws.Range("L74").Value = 42474.5
$ws.Range("N74").Value = -44222.5
$ws.Range("H77").Value = 42474.5
$ws.Range("J77").Value = 42474.5
$ws.Range("L77").Value = 127423.5
$ws.Range("N77").Value = -136159.5
$ws.Range("H113").Value = 2710.4666
$ws.Range("I113").Value = 2606.3076
$ws.Range("K113").Value = 2606.3076
$ws.Range("M113").Value = -436.3076000000001
$ws.Range("H134").Value = 14149.5
$ws.Range("I134").Value = 13024.5
$ws.Range("K134").Value = 39073.5
$ws.Range("M134").Value = -36538.5
$ws.Range("H136").Value = 5176.5386
$ws.Range("I136").Value = 2548.75
$ws.Range("K136").Value = 7646.25
$ws.Range("M136").Value = -5096.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 137.81818
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 96.545456
$ws.Range("I2").Value = 34.8
$ws.Range("K2").Value = 34.8
$ws.Range("M2").Value = 78.2
$ws.Range("H10").Value = 9144.833000000001
$ws.Range("I10").Value = 6467.5
$ws.Range("K10").Value = 6467.5
$ws.Range("M10").Value = -6298.5
$ws.Range("H97").Value = 563.8
$ws.Range("I97").Value = 469.0303
$ws.Range("J97").Value = 1010.5714
$ws.Range("K97").Value = 469.0303
$ws.Range("L97").Value = 1010.5714
$ws.Range("M97").Value = 26.96969999999999
$ws.Range("N97").Value = -2002.5714
$ws.Range("H132").Value = 4265.625
$ws.Range("I132").Value = 3354.3333
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 10062.9999
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -7532.999899999999
$ws.Range("N132").Value = -26058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 730.9666999999999
$ws.Range("I46").Value = 686.94446
$ws.Range("J46").Value = 797
$ws.Range("K46").Value = 686.94446
$ws.Range("L46").Value = 797
$ws.Range("M46").Value = -498.94446
$ws.Range("N46").Value = -1173
$ws.Range("H61").Value = 1888.4
$ws.Range("I61").Value = 1804.3846
$ws.Range("K61").Value = 1804.3846
$ws.Range("M61").Value = -1602.3846
$ws.Range("H82").Value = 1544.75
$ws.Range("I82").Value = 1966.5
$ws.Range("K82").Value = 1966.5
$ws.Range("M82").Value = -1605.5
$ws.Range("H85").Value = 1544.75
$ws.Range("I85").Value = 1966.5
$ws.Range("K85").Value = 1966.5
$ws.Range("M85").Value = -718.5
$ws.Range("H93").Value = 1861.875
$ws.Range("I93").Value = 1970.6666
$ws.Range("K93").Value = 1970.6666
$ws.Range("M93").Value = -722.6666
$ws.Range("H100").Value = 860682.3
$ws.Range("I100").Value = 7177
$ws.Range("K100").Value = 7177
$ws.Range("M100").Value = -6636
$ws.Range("H113").Value = 1888.4
$ws.Range("I113").Value = 1804.3846
$ws.Range("K113").Value = 1804.3846
$ws.Range("M113").Value = 365.6153999999999
$ws.Range("H132").Value = 6946913
$ws.Range("I132").Value = 6946913
$ws.Range("K132").Value = 20840739
$ws.Range("M132").Value = -20838209
$ws.Range("H136").Value = 5557690
$ws.Range("I136").Value = 6538100.5
$ws.Range("J136").Value = 2028
$ws.Range("K136").Value = 19614301.5
$ws.Range("L136").Value = 6084
$ws.Range("M136").Value = -19611751.5
$ws.Range("N136").Value = -11184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 427.85715
$ws.Range("I113").Value = 311.30768
$ws.Range("K113").Value = 933.92304
$ws.Range("M113").Value = 1236.07696
$ws.Range("H126").Value = 8115.091
$ws.Range("J126").Value = 6247.5
$ws.Range("L126").Value = 18742.5
$ws.Range("N126").Value = -23682.5
